$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-06-03 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-04 Wednesday", 2)

# Update the division problems in the single table. Rows 1, 5, 9, 13, 17
# contain the problems (the rows in between are blank answer rows).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "26÷9="
$t.Cell(1, 2).Range.Text = "15÷8="
$t.Cell(1, 3).Range.Text = "60÷6="
$t.Cell(1, 4).Range.Text = "42÷8="
$t.Cell(1, 5).Range.Text = "53÷9="

$t.Cell(5, 1).Range.Text = "44÷7="
$t.Cell(5, 2).Range.Text = "44÷3="
$t.Cell(5, 3).Range.Text = "88÷3="
$t.Cell(5, 4).Range.Text = "71÷3="
$t.Cell(5, 5).Range.Text = "94÷3="

$t.Cell(9, 1).Range.Text = "52÷2="
$t.Cell(9, 2).Range.Text = "13÷4="
$t.Cell(9, 3).Range.Text = "87÷2="
$t.Cell(9, 4).Range.Text = "67÷3="
$t.Cell(9, 5).Range.Text = "59÷3="

$t.Cell(13, 1).Range.Text = "21÷9="
$t.Cell(13, 2).Range.Text = "78÷5="
$t.Cell(13, 3).Range.Text = "57÷6="
$t.Cell(13, 4).Range.Text = "93÷9="
$t.Cell(13, 5).Range.Text = "86÷7="

$t.Cell(17, 1).Range.Text = "26÷9="
$t.Cell(17, 2).Range.Text = "64÷8="
$t.Cell(17, 3).Range.Text = "52÷5="
$t.Cell(17, 4).Range.Text = "62÷3="
$t.Cell(17, 5).Range.Text = "60÷6="
